# Update the QR-code text box ("Meldungs-ID" / {%DemisIdQrImage}) in the
# report header: wrap the existing QR-image placeholder with a
# "{#DemisIdQrImage} ... {%DemisIdQrImage} ... {/DemisIdQrImage}" template
# block (adds an opening comment tag before the image placeholder and a
# closing tag after the "Meldungs-ID" caption), all inside the legacy VML
# text box (w:pict/v:shape) that lives in the report-title paragraph.

$d = $word.ActiveDocument

# Locate the paragraph that owns the floating text box (the VML <w:pict>
# is not reachable through Range.Text / Find, so we find it by scanning
# the paragraphs' raw OOXML for the marker we know is inside the shape).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $para = $d.Paragraphs.Item($i)
    $openXml = $para.Range.WordOpenXML
    if ($openXml.Contains("DemisIdQrImage")) {
        $target = $para
        break
    }
}

if ($target -eq $null) {
    throw "Could not find paragraph containing the DemisIdQrImage text box"
}

$range = $target.Range
$xml = $range.WordOpenXML

# --- 1. anchorId / shape id refresh (cosmetic ids Word re-mints on edit) ---
$xml = $xml.Replace('w14:anchorId="4060B227"', 'w14:anchorId="6B3999DC"')
$xml = $xml.Replace('id="_x0000_s2050"', 'id="_x0000_s2052"')

# --- 2. Insert the "{#" "DemisIdQrImage" "}" runs right before the
#        existing "{%DemisIdQrImage}" run in the first text-box paragraph.
$oldImageRun = '<w:r w:rsidRPr="006F7F3C"><w:rPr><w:sz w:val="16"/><w:szCs w:val="16"/></w:rPr><w:t>{%DemisIdQrImage}</w:t></w:r>'
if (-not $xml.Contains($oldImageRun)) {
    throw "Could not find the {%DemisIdQrImage} run to anchor the edit"
}

$newRunPropsOpen = '<w:rPr><w:rFonts w:ascii="Arial" w:hAnsi="Arial" w:cs="Arial"/><w:sz w:val="12"/><w:szCs w:val="12"/><w:lang w:val="en-US"/></w:rPr>'
$openTagRuns = '<w:r>' + $newRunPropsOpen + '<w:t>{#</w:t></w:r>' + `
               '<w:r>' + $newRunPropsOpen + '<w:t>DemisIdQrImage</w:t></w:r>' + `
               '<w:r>' + $newRunPropsOpen + '<w:t>}</w:t></w:r>'

$xml = $xml.Replace($oldImageRun, $openTagRuns + $oldImageRun)

# --- 3. Insert the "{" "/" "DemisIdQrImage" "}" runs right after the
#        "Meldungs-ID" run in the second text-box paragraph.
$oldLabelRun = '<w:r w:rsidRPr="000E6D63"><w:rPr><w:sz w:val="12"/><w:szCs w:val="12"/></w:rPr><w:t>Meldungs-ID</w:t></w:r>'
if (-not $xml.Contains($oldLabelRun)) {
    throw "Could not find the Meldungs-ID run to anchor the edit"
}

$closeTagRuns = '<w:r>' + $newRunPropsOpen + '<w:t>{</w:t></w:r>' + `
                '<w:r>' + $newRunPropsOpen + '<w:t>/</w:t></w:r>' + `
                '<w:r>' + $newRunPropsOpen + '<w:t>DemisIdQrImage</w:t></w:r>' + `
                '<w:r>' + $newRunPropsOpen + '<w:t>}</w:t></w:r>'

$xml = $xml.Replace($oldLabelRun, $oldLabelRun + $closeTagRuns)

# Push the rebuilt OOXML back into the document, replacing the paragraph's
# content in place.
[void]$range.InsertXML($xml)
